$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the review comments on B2/J2/O2 (no longer needed)
$ws.Range("B2").Comment.Delete()
$ws.Range("J2").Comment.Delete()
$ws.Range("O2").Comment.Delete()

# Drop the now-unused "Attrib_Cond"/"Val_Cond" columns (F:G) and insert a
# fresh single column at G for the new "Top_Check" qualifier
$ws.Range("F1:G1").EntireColumn.Delete()
$ws.Range("G1").EntireColumn.Insert()

# Header row: rename the old "CUR" column (now shifted into F) and label
# the newly inserted column
$ws.Range("F3").Value = "CURR"
$ws.Range("G3").Value = "Top_Check"

# Data rows: populate the new Top_Check column
$ws.Range("G4").Value = "O"
$ws.Range("G5").Value = "O"
$ws.Range("G6").Value = "O"

# Widen the technology-qualifier filter to also exclude navigation techs
# (leading apostrophe preserves the cell's existing quote-prefix formatting)
$ws.Range("J4").Value = "'*,-T-A*INT*,-T-NAV*"

# Match the author's final cursor position
$ws.Range("F4").Select()
